# Auto-generated update of Sheets (Phantom_Profits) market-price refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 291.52173
$ws.Range("I33").Value = 276.4737
$ws.Range("J33").Value = 363
$ws.Range("K33").Value = 276.4737
$ws.Range("L33").Value = 363
$ws.Range("M33").Value = -47.47370000000001
$ws.Range("N33").Value = -821

# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 3308.9092
$ws.Range("I40").Value = 1561.3077
$ws.Range("K40").Value = 1561.3077
$ws.Range("M40").Value = -1386.3077

# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 3196.4
$ws.Range("J43").Value = 3327
$ws.Range("L43").Value = 3327
$ws.Range("N43").Value = -3465

# Row 55 (Leve Item ID 5517)
$ws.Range("H55").Value = 773.5333000000001
$ws.Range("I55").Value = 267.55554
$ws.Range("K55").Value = 267.55554
$ws.Range("M55").Value = -53.55554000000001

# Row 88 (Leve Item ID 12608)
$ws.Range("H88").Value = 2881.5
$ws.Range("J88").Value = 2894.5
$ws.Range("L88").Value = 2894.5
$ws.Range("N88").Value = -3706.5

# Row 91 (Leve Item ID 12608)
$ws.Range("H91").Value = 2881.5
$ws.Range("J91").Value = 2894.5
$ws.Range("L91").Value = 2894.5
$ws.Range("N91").Value = -5702.5

$ws = $wb.Worksheets.Item("ARM")
# Row 4 (Leve Item ID 5071)
$ws.Range("H4").Value = 566.8570999999999
$ws.Range("I4").Value = 713.6
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 713.6
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -597.6
$ws.Range("N4").Value = -432

# Row 11 (Leve Item ID 3767)
$ws.Range("H11").Value = 3984.3333
$ws.Range("J11").Value = 975
$ws.Range("L11").Value = 975
$ws.Range("N11").Value = -1263

# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 8777.700000000001
$ws.Range("I32").Value = 8187.0527
$ws.Range("K32").Value = 8187.0527
$ws.Range("M32").Value = -7900.0527

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 9621.666999999999
$ws.Range("I61").Value = 9621.666999999999
$ws.Range("K61").Value = 9621.666999999999
$ws.Range("M61").Value = -9409.666999999999

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2721.0715
$ws.Range("I122").Value = 3961.875
$ws.Range("K122").Value = 11885.625
$ws.Range("M122").Value = -9435.625

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3397.6206
$ws.Range("I132").Value = 3427.0454
$ws.Range("K132").Value = 10281.1362
$ws.Range("M132").Value = -7751.136200000001

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 9621.666999999999
$ws.Range("I136").Value = 9621.666999999999
$ws.Range("K136").Value = 28865.001
$ws.Range("M136").Value = -26315.001

$ws = $wb.Worksheets.Item("BSM")
# Row 35 (Leve Item ID 2350)
$ws.Range("H35").Value = 60073.5
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 60073.5
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 60073.5
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -60693.5

# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 1618.9131
$ws.Range("I94").Value = 4161.6665
$ws.Range("K94").Value = 4161.6665
$ws.Range("M94").Value = -3710.6665

# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 1782.8462
$ws.Range("J99").Value = 3824.25
$ws.Range("L99").Value = 3824.25
$ws.Range("N99").Value = -6820.25

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 331.27274
$ws.Range("I7").Value = 130.5
$ws.Range("J7").Value = 866.6667
$ws.Range("K7").Value = 130.5
$ws.Range("L7").Value = 866.6667
$ws.Range("M7").Value = -17.5
$ws.Range("N7").Value = -1092.6667

# Row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 5716671.5
$ws.Range("I22").Value = 2234
$ws.Range("J22").Value = 10002500
$ws.Range("K22").Value = 2234
$ws.Range("L22").Value = 10002500
$ws.Range("M22").Value = -1884
$ws.Range("N22").Value = -10003200

# Row 59 (Leve Item ID 1942)
$ws.Range("H59").Value = 70000
$ws.Range("I59").Value = 70000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 70000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -68855
$ws.Range("N59").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 640.6
$ws.Range("J5").Value = 597.6667
$ws.Range("L5").Value = 1793.0001
$ws.Range("N5").Value = -2017.0001

# Row 8 (Leve Item ID 16734)
$ws.Range("H8").Value = 130
$ws.Range("I8").Value = 130
$ws.Range("K8").Value = 390
$ws.Range("M8").Value = -251

# Row 12 (Leve Item ID 4854)
$ws.Range("H12").Value = 43.94737
$ws.Range("I12").Value = 29.8
$ws.Range("J12").Value = 49
$ws.Range("K12").Value = 89.40000000000001
$ws.Range("L12").Value = 147
$ws.Range("M12").Value = 83.59999999999999
$ws.Range("N12").Value = -493

# Row 29 (Leve Item ID 4698)
$ws.Range("H29").Value = 375
$ws.Range("J29").Value = 375
$ws.Range("L29").Value = 1125
$ws.Range("N29").Value = -1679

# Row 32 (Leve Item ID 4731)
$ws.Range("H32").Value = 1497
$ws.Range("J32").Value = 999
$ws.Range("L32").Value = 2997
$ws.Range("N32").Value = -3563

# Row 46 (Leve Item ID 4701)
$ws.Range("H46").Value = 2576.6
$ws.Range("I46").Value = 2995.75
$ws.Range("J46").Value = 900
$ws.Range("K46").Value = 8987.25
$ws.Range("L46").Value = 2700
$ws.Range("M46").Value = -8896.25
$ws.Range("N46").Value = -2882

# Row 61 (Leve Item ID 4727)
$ws.Range("H61").Value = 290.375
$ws.Range("I61").Value = 189
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 567
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -352
$ws.Range("N61").Value = -3430

# Row 93 (Leve Item ID 19808)
$ws.Range("H93").Value = 899
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 999
$ws.Range("I132").Value = 999
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8991
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6461
$ws.Range("N132").ClearContents()

# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 640.6
$ws.Range("J135").Value = 597.6667
$ws.Range("L135").Value = 5379.0003
$ws.Range("N135").Value = -10449.0003

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 13.333333
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 17.5
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 17.5
$ws.Range("M2").Value = 108
$ws.Range("N2").Value = -243.5

# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 10166.667
$ws.Range("I70").Value = 10500
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 10500
$ws.Range("L70").Value = 10000
$ws.Range("M70").Value = -10230
$ws.Range("N70").Value = -10540

# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 10166.667
$ws.Range("I73").Value = 10500
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 10500
$ws.Range("L73").Value = 10000
$ws.Range("M73").Value = -9564
$ws.Range("N73").Value = -11872

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 4174.75
$ws.Range("I22").Value = 899
$ws.Range("J22").Value = 5266.6665
$ws.Range("K22").Value = 899
$ws.Range("L22").Value = 5266.6665
$ws.Range("M22").Value = -604
$ws.Range("N22").Value = -5856.6665

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 4174.75
$ws.Range("I27").Value = 899
$ws.Range("J27").Value = 5266.6665
$ws.Range("K27").Value = 899
$ws.Range("L27").Value = 5266.6665
$ws.Range("M27").Value = -792
$ws.Range("N27").Value = -5480.6665

# Row 33 (Leve Item ID 4106)
$ws.Range("H33").Value = 11111
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 11111
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 11111
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -11691

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 1483
$ws.Range("I46").Value = 1102.4375
$ws.Range("K46").Value = 1102.4375
$ws.Range("M46").Value = -914.4375

# Row 116 (Leve Item ID 26133)
$ws.Range("H116").Value = 21748.375
$ws.Range("J116").Value = 19998.143
$ws.Range("L116").Value = 19998.143
$ws.Range("N116").Value = -29176.143

$ws = $wb.Worksheets.Item("WVR")
# Row 12 (Leve Item ID 3316)
$ws.Range("H12").Value = 10500
$ws.Range("J12").Value = 10500
$ws.Range("L12").Value = 10500
$ws.Range("N12").Value = -10784

# Row 31 (Leve Item ID 3052)
$ws.Range("H31").Value = 4000
$ws.Range("I31").Value = 4000
$ws.Range("K31").Value = 4000
$ws.Range("M31").Value = -3652

# Row 54 (Leve Item ID 3413)
$ws.Range("H54").Value = 17332.666
$ws.Range("J54").Value = 31998
$ws.Range("L54").Value = 31998
$ws.Range("N54").Value = -33038

# Row 96 (Leve Item ID 19977)
$ws.Range("H96").Value = 5578.6665
$ws.Range("I96").Value = 5344.4
$ws.Range("K96").Value = 5344.4
$ws.Range("M96").Value = -3971.4

# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 397.33334
$ws.Range("I113").Value = 231.2
$ws.Range("K113").Value = 693.5999999999999
$ws.Range("M113").Value = 1476.4

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 3637
$ws.Range("I132").Value = 3000.875
$ws.Range("K132").Value = 9002.625
$ws.Range("M132").Value = -6472.625
